# Generate Report for Handoff
#
# Updates the localization-status report:
#  - "Overview" sheet: Latest HO Xliff Generate Date (col G) for the rows
#    that were handed off at 2016-08-24 04:20:12 -> 04:20:35
#  - "de-de" sheet: Latest Handoff Datetime (col H) same timestamp bump
#  - "zh-cn" sheet: Latest Handoff Datetime (col H) 04:20:00 -> 04:20:30
#  - "zh-cn" / "de-de" sheets: Priority (col E) set to "ht" for the same
#    rows (was blank)

$wb = $excel.ActiveWorkbook

$rows = @(7, 9, 10, 11, 12, 14)

$wsOverview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $wsOverview.Cells.Item($r, 7).Value = "2016-08-24 04:20:35"
}

$wsZhCn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $wsZhCn.Cells.Item($r, 8).Value = "2016-08-24 04:20:30"
    $wsZhCn.Cells.Item($r, 5).Value = "ht"
}

$wsDeDe = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $wsDeDe.Cells.Item($r, 8).Value = "2016-08-24 04:20:35"
    $wsDeDe.Cells.Item($r, 5).Value = "ht"
}
